$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 227.282303
$ws.Range("H2").Value = 681.846909
$ws.Range("I2").Value = 0.6094595465130797
$ws.Range("J2").Value = 0.6094595465130795
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 2881.757170111459
$ws.Range("R2").Value = 25935.81453100313
$ws.Range("S2").Value = 0.6016228807848414
$ws.Range("T2").Value = 0.6016228807848414

# Row 3
$ws.Range("G3").Value = 227.282303
$ws.Range("H3").Value = 681.846909
$ws.Range("I3").Value = 0.6094595465130797
$ws.Range("J3").Value = 0.6094595465130795
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 37.53741483810633
$ws.Range("R3").Value = 337.836733542957
$ws.Range("S3").Value = 0.007836665728238205
$ws.Range("T3").Value = 0.007836665728238204

# Row 4
$ws.Range("I4").Value = 0.3727053955221387
$ws.Range("J4").Value = 0.3727053955221385
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 1762.293251504758
$ws.Range("R4").Value = 15860.63926354282
$ws.Range("S4").Value = 0.3679130059098527
$ws.Range("T4").Value = 0.3679130059098526

# Row 5
$ws.Range("I5").Value = 0.3727053955221387
$ws.Range("J5").Value = 0.3727053955221385
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("Q5").Value = 22.95541537442267
$ws.Range("S5").Value = 0.004792389612285983
$ws.Range("T5").Value = 0.004792389612285982

# Row 6
$ws.Range("G6").Value = 6.377814666666666
$ws.Range("H6").Value = 19.133444
$ws.Range("I6").Value = 0.0171021675827138
$ws.Range("J6").Value = 0.01710216758271379
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 80.86557071409422
$ws.Range("R6").Value = 727.790136426848
$ws.Range("S6").Value = 0.01688226132094329
$ws.Range("T6").Value = 0.01688226132094329

# Row 7
$ws.Range("G7").Value = 6.377814666666666
$ws.Range("H7").Value = 19.133444
$ws.Range("I7").Value = 0.0171021675827138
$ws.Range("J7").Value = 0.01710216758271379
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 1.053344988779111
$ws.Range("R7").Value = 9.480104899011998
$ws.Range("S7").Value = 0.0002199062617705068
$ws.Range("T7").Value = 0.0002199062617705067

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2733126666666666
$ws.Range("H8").Value = 0.8199379999999999
$ws.Range("I8").Value = 0.0007328903820679218
$ws.Range("J8").Value = 0.0007328903820679217
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 3.465385234366222
$ws.Range("R8").Value = 31.188467109296
$ws.Range("S8").Value = 0.0007234665950871991
$ws.Range("T8").Value = 0.0007234665950871991

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2733126666666666
$ws.Range("H9").Value = 0.8199379999999999
$ws.Range("I9").Value = 0.0007328903820679218
$ws.Range("J9").Value = 0.0007328903820679217
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 0.04513968229711111
$ws.Range("R9").Value = 0.406257140674
$ws.Range("S9").Value = 0.000009423786980722644
$ws.Range("T9").Value = 0.000009423786980722644
